# Updates cryptos list prices/volumes (and a couple of row-content swaps)
# to match the latest scrape, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.746.30"
$ws.Range("E2").Value = "  -0.34%  "

$ws.Range("D3").Value = "1.633.56"
$ws.Range("E3").Value = "  -0.49%  "

$ws.Range("D5").Value = "'215.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("D6").Value = "'0.504"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.61%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("E9").Value = "  -1.22%  "

$ws.Range("D10").Value = "'19.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.25%  "

$ws.Range("D11").Value = "'0.0785"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.78%  "

$ws.Range("E12").Value = "  -0.63%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.859.26"
$ws.Range("E13").Value = "  -0.42%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.633.94"
$ws.Range("E14").Value = "  -0.80%  "

$ws.Range("E15").Value = "  -1.18%  "

$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("E16").Value = "  -0.12%  "

$ws.Range("D17").Value = "'62.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.14%  "

$ws.Range("D18").Value = "25.764.47"
$ws.Range("E18").Value = "  -0.32%  "

$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("E20").Value = "  +0.99%  "

$ws.Range("D21").Value = "'193.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.50%  "

$ws.Range("D22").Value = "'9.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").Value = "'6.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.17%  "

$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("D25").Value = "'1.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.17%  "

$ws.Range("D26").Value = "'140.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.60%  "

$ws.Range("D27").Value = "'0.122"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.21%  "

$ws.Range("D28").Value = "'6.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.82%  "

$ws.Range("E29").Value = "  -0.35%  "

$ws.Range("E30").Value = "  -0.32%  "

$ws.Range("D31").Value = "'0.0493"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.30%  "

$ws.Range("E32").Value = "  +1.10%  "

$ws.Range("E33").Value = "  -0.37%  "

$ws.Range("D34").Value = "'1.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.00%  "

$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("E36").Value = "  -0.57%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.548"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.35%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.120.04"
$ws.Range("E38").Value = "  -1.23%  "

$ws.Range("D39").Value = "'2.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.78%  "

$ws.Range("E40").Value = "  -1.05%  "

$ws.Range("E41").Value = "  +0.56%  "

$ws.Range("E42").Value = "  +1.94%  "

$ws.Range("D43").Value = "'99.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.53%  "

$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").Value = "1.768.36"
$ws.Range("E45").Value = "  -0.51%  "

$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("E46").Value = "  -1.84%  "

$ws.Range("D47").Value = "'55.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.02%  "

$ws.Range("E48").Value = "  -2.41%  "

$ws.Range("D50").Value = "'7.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.31%  "

$ws.Range("D51").Value = "'2.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.90%  "
